$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$arr2 = New-Object 'object[,]' 1,23
$arr2[0,0] = 1.81399941444397
$arr2[0,1] = 1
$arr2[0,2] = 6402.475950800274
$arr2[0,3] = 0.2435465729302398
$arr2[0,4] = 0.2052815502316371
$arr2[0,5] = 0.178694355199945
$arr2[0,6] = 0.1683892237606172
$arr2[0,7] = 0.1578628390298868
$arr2[0,8] = 0.1509940945340159
$arr2[0,9] = 0.1425904236156144
$arr2[0,10] = 0.1354602983680519
$arr2[0,11] = 0.1332350293336609
$arr2[0,12] = 0.1328878184165517
$arr2[0,13] = 0.1328878184165517
$arr2[0,14] = 0.1328878184165517
$arr2[0,15] = 0.1328878184165517
$arr2[0,16] = 0.1328878184165517
$arr2[0,17] = 0.1328045994308045
$arr2[0,18] = 0.1328045994308045
$arr2[0,19] = 0.1328045994308045
$arr2[0,20] = 0.1328045994308045
$arr2[0,21] = 0.1328045994308045
$arr2[0,22] = 0.1328045994308045
$ws.Range("C2:Y2").Value2 = $arr2

$arr3 = New-Object 'object[,]' 1,23
$arr3[0,0] = 1.485998868942261
$arr3[0,1] = 2
$arr3[0,2] = 6428.228550823073
$arr3[0,3] = 0.2435465729302398
$arr3[0,4] = 0.2052815502316371
$arr3[0,5] = 0.178694355199945
$arr3[0,6] = 0.1683892237606172
$arr3[0,7] = 0.1634031916298378
$arr3[0,8] = 0.1565002134468439
$arr3[0,9] = 0.1521304853218897
$arr3[0,10] = 0.1469182874437504
$arr3[0,11] = 0.1469182874437504
$arr3[0,12] = 0.1459895682571866
$arr3[0,13] = 0.1427917377320834
$arr3[0,14] = 0.1416120370577794
$arr3[0,15] = 0.1416120370577794
$arr3[0,16] = 0.1415380397343942
$arr3[0,17] = 0.1415380397343942
$arr3[0,18] = 0.1415380397343942
$arr3[0,19] = 0.1414888013901989
$arr3[0,20] = 0.1414888013901989
$arr3[0,21] = 0.1414888013901989
$arr3[0,22] = 0.141306599431249
$ws.Range("C3:Y3").Value2 = $arr3

$arr4 = New-Object 'object[,]' 1,23
$arr4[0,0] = 1.717000961303711
$arr4[0,1] = 1
$arr4[0,2] = 6710.309206785329
$arr4[0,3] = 0.2435465729302398
$arr4[0,4] = 0.2052815502316371
$arr4[0,5] = 0.178694355199945
$arr4[0,6] = 0.1683892237606172
$arr4[0,7] = 0.1549879400411887
$arr4[0,8] = 0.1487521583260372
$arr4[0,9] = 0.147463193855715
$arr4[0,10] = 0.1451535800088008
$arr4[0,11] = 0.1388052476956204
$arr4[0,12] = 0.1388052476956204
$arr4[0,13] = 0.1388052476956204
$arr4[0,14] = 0.1388052476956204
$arr4[0,15] = 0.1388052476956204
$arr4[0,16] = 0.1388052476956204
$arr4[0,17] = 0.1388052476956204
$arr4[0,18] = 0.1388052476956204
$arr4[0,19] = 0.1388052476956204
$arr4[0,20] = 0.1388052476956204
$arr4[0,21] = 0.1388052476956204
$arr4[0,22] = 0.1388052476956204
$ws.Range("C4:Y4").Value2 = $arr4

$arr5 = New-Object 'object[,]' 1,23
$arr5[0,0] = 1.679980993270874
$arr5[0,1] = 2
$arr5[0,2] = 6401.079491542452
$arr5[0,3] = 0.2435465729302398
$arr5[0,4] = 0.2052815502316371
$arr5[0,5] = 0.178694355199945
$arr5[0,6] = 0.1683892237606172
$arr5[0,7] = 0.149634652393159
$arr5[0,8] = 0.1464057319029364
$arr5[0,9] = 0.1447392307779213
$arr5[0,10] = 0.1413962662297761
$arr5[0,11] = 0.1413962662297761
$arr5[0,12] = 0.1409155379320654
$arr5[0,13] = 0.1408344258264169
$arr5[0,14] = 0.1408344258264169
$arr5[0,15] = 0.1408344258264169
$arr5[0,16] = 0.1408344258264169
$arr5[0,17] = 0.1407773780027768
$arr5[0,18] = 0.1407773780027768
$arr5[0,19] = 0.1407773780027768
$arr5[0,20] = 0.1407773780027768
$arr5[0,21] = 0.1407773780027768
$arr5[0,22] = 0.1407773780027768
$ws.Range("C5:Y5").Value2 = $arr5

$ws.Range("C6").Value2 = 1.536943674087524
$arr6 = New-Object 'object[,]' 1,21
$arr6[0,0] = 6411.447511061301
$arr6[0,1] = 0.2435465729302398
$arr6[0,2] = 0.2052815502316371
$arr6[0,3] = 0.178694355199945
$arr6[0,4] = 0.1683892237606172
$arr6[0,5] = 0.1528855784489951
$arr6[0,6] = 0.1528855784489951
$arr6[0,7] = 0.1481130607365902
$arr6[0,8] = 0.140216055188859
$arr6[0,9] = 0.1375117771559863
$arr6[0,10] = 0.1352056187549749
$arr6[0,11] = 0.1343482026507863
$arr6[0,12] = 0.1339639066960592
$arr6[0,13] = 0.1339639066960592
$arr6[0,14] = 0.1339413678551649
$arr6[0,15] = 0.1329794836464191
$arr6[0,16] = 0.1329794836464191
$arr6[0,17] = 0.1329794836464191
$arr6[0,18] = 0.1329794836464191
$arr6[0,19] = 0.1329794836464191
$arr6[0,20] = 0.1329794836464191
$ws.Range("E6:Y6").Value2 = $arr6

$arr7 = New-Object 'object[,]' 1,23
$arr7[0,0] = 1.613997220993042
$arr7[0,1] = 1
$arr7[0,2] = 6398.337455565823
$arr7[0,3] = 0.2435465729302398
$arr7[0,4] = 0.2052815502316371
$arr7[0,5] = 0.178694355199945
$arr7[0,6] = 0.1683892237606172
$arr7[0,7] = 0.1529483792876766
$arr7[0,8] = 0.1465864945995281
$arr7[0,9] = 0.1400525914119534
$arr7[0,10] = 0.1341628036367663
$arr7[0,11] = 0.1331957942852961
$arr7[0,12] = 0.1328847939642564
$arr7[0,13] = 0.1327323239271318
$arr7[0,14] = 0.1327323239271318
$arr7[0,15] = 0.1327323239271318
$arr7[0,16] = 0.1327287366441272
$arr7[0,17] = 0.1327287366441272
$arr7[0,18] = 0.1327287366441272
$arr7[0,19] = 0.1327239270090803
$arr7[0,20] = 0.1327239270090803
$arr7[0,21] = 0.1327239270090803
$arr7[0,22] = 0.1327239270090803
$ws.Range("C7:Y7").Value2 = $arr7

$ws.Range("C8").Value2 = 1.538951635360718
$arr8 = New-Object 'object[,]' 1,21
$arr8[0,0] = 6404.544272658598
$arr8[0,1] = 0.2435465729302398
$arr8[0,2] = 0.2052815502316371
$arr8[0,3] = 0.178694355199945
$arr8[0,4] = 0.1683892237606172
$arr8[0,5] = 0.1635707437378576
$arr8[0,6] = 0.1489498349514058
$arr8[0,7] = 0.1419932390705304
$arr8[0,8] = 0.1395383037890676
$arr8[0,9] = 0.1352748415216457
$arr8[0,10] = 0.1329749553030454
$arr8[0,11] = 0.1329749553030454
$arr8[0,12] = 0.1329749553030454
$arr8[0,13] = 0.1329233210192076
$arr8[0,14] = 0.1329233210192076
$arr8[0,15] = 0.1328566499580744
$arr8[0,16] = 0.1328566499580744
$arr8[0,17] = 0.1328566499580744
$arr8[0,18] = 0.1328541049716566
$arr8[0,19] = 0.1328541049716566
$arr8[0,20] = 0.1328449175956841
$ws.Range("E8:Y8").Value2 = $arr8

$arr9 = New-Object 'object[,]' 1,23
$arr9[0,0] = 1.568004369735718
$arr9[0,1] = 2
$arr9[0,2] = 6403.751574711064
$arr9[0,3] = 0.2435465729302398
$arr9[0,4] = 0.2052815502316371
$arr9[0,5] = 0.178694355199945
$arr9[0,6] = 0.1669630682704453
$arr9[0,7] = 0.1529917639035948
$arr9[0,8] = 0.1526767706537398
$arr9[0,9] = 0.1463237862372542
$arr9[0,10] = 0.141050699889859
$arr9[0,11] = 0.141050699889859
$arr9[0,12] = 0.1409017221422948
$arr9[0,13] = 0.1409017221422948
$arr9[0,14] = 0.1408467084297077
$arr9[0,15] = 0.1408467084297077
$arr9[0,16] = 0.1408467084297077
$arr9[0,17] = 0.1408467084297077
$arr9[0,18] = 0.1408467084297077
$arr9[0,19] = 0.1408467084297077
$arr9[0,20] = 0.1408467084297077
$arr9[0,21] = 0.1408382574357522
$arr9[0,22] = 0.1408294653939778
$ws.Range("C9:Y9").Value2 = $arr9

$arr10 = New-Object 'object[,]' 1,23
$arr10[0,0] = 1.589019060134888
$arr10[0,1] = 1
$arr10[0,2] = 6401.560080297098
$arr10[0,3] = 0.2435465729302398
$arr10[0,4] = 0.2052815502316371
$arr10[0,5] = 0.178694355199945
$arr10[0,6] = 0.1683892237606172
$arr10[0,7] = 0.1567560626070753
$arr10[0,8] = 0.1530316141783553
$arr10[0,9] = 0.1457413466862789
$arr10[0,10] = 0.1392755709879787
$arr10[0,11] = 0.1363704599490802
$arr10[0,12] = 0.1357716411374943
$arr10[0,13] = 0.1333464916588328
$arr10[0,14] = 0.132885143174724
$arr10[0,15] = 0.132885143174724
$arr10[0,16] = 0.132885143174724
$arr10[0,17] = 0.132885143174724
$arr10[0,18] = 0.132885143174724
$arr10[0,19] = 0.1328228891261538
$arr10[0,20] = 0.1328228891261538
$arr10[0,21] = 0.1327867462046218
$arr10[0,22] = 0.1327867462046218
$ws.Range("C10:Y10").Value2 = $arr10

$arr11 = New-Object 'object[,]' 1,23
$arr11[0,0] = 1.57801079750061
$arr11[0,1] = 2
$arr11[0,2] = 6398.536599706097
$arr11[0,3] = 0.2435465729302398
$arr11[0,4] = 0.2052815502316371
$arr11[0,5] = 0.178694355199945
$arr11[0,6] = 0.1683892237606172
$arr11[0,7] = 0.1549618883998216
$arr11[0,8] = 0.14869410273477
$arr11[0,9] = 0.1448079053246598
$arr11[0,10] = 0.1448079053246598
$arr11[0,11] = 0.1440024987646787
$arr11[0,12] = 0.1410954357640267
$arr11[0,13] = 0.1407601490569694
$arr11[0,14] = 0.1407445688223459
$arr11[0,15] = 0.1407445688223459
$arr11[0,16] = 0.1407278089611325
$arr11[0,17] = 0.1407278089611325
$arr11[0,18] = 0.1407278089611325
$arr11[0,19] = 0.1407278089611325
$arr11[0,20] = 0.1407278089611325
$arr11[0,21] = 0.1407278089611325
$arr11[0,22] = 0.1407278089611325
$ws.Range("C11:Y11").Value2 = $arr11
